$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.474.81'

$ws.Range("D3").Value = '1.824.27'
$ws.Range("E3").Value = '  +1.66%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '316.90'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.03%  '

$ws.Range("E6").Value = '  +0.03%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5410'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.79%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.4034'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +6.98%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07648'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.47%  '

$ws.Range("E10").Value = '  +2.42%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '41.91'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.61%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '6.334'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.73%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '7.651'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +6.08%  '

$ws.Range("E14").Value = '  +0.08%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '20.91'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.40%  '

$ws.Range("D16").Value = '1.823.92'
$ws.Range("E16").Value = '  +2.11%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001090'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +3.13%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '89.82'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.74%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06610'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.25%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.70'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("E21").Value = '  +0.06%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.067'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.86%  '

$ws.Range("D23").Value = '28.477.99'
$ws.Range("E23").Value = '  +1.24%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '11.11'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.67%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.276'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +8.59%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '157.94'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +2.01%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.456'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +7.66%  '

$ws.Range("D29").Value = '2.034.15'
$ws.Range("E29").Value = '  +2.07%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '123.98'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +2.78%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.123'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.59%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.679'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.17%  '

$ws.Range("E34").Value = '  -0.37%  '

$ws.Range("E35").Value = '  +12.21%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.2238'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.93%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.02345'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.61%  '

$ws.Range("E38").Value = '  +3.82%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.855'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +4.84%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.6294'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.03%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '11.35'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +2.43%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.188'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.49%  '

$ws.Range("E43").Value = '  +0.06%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.400'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.23%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '13.44'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.98%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.698'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5857'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.33%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '125.41'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.002'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.99%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.197'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.64%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06878'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.96%  '
